# Insert a new data row at row 500 (shifting existing rows 500-530 down to 501-531)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(500).Insert()

$ws.Cells.Item(500, 1).Value = 3
$ws.Cells.Item(500, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(500, 3).Value = "Coquimbo"
$ws.Cells.Item(500, 4).Value = 44826
$ws.Cells.Item(500, 5).Value = 5
$ws.Cells.Item(500, 6).Value = 100112032
$ws.Cells.Item(500, 7).Value = "Zapallo italiano"
$ws.Cells.Item(500, 8).Value = "Sin especificar"
$ws.Cells.Item(500, 9).Value = "Primera"
$ws.Cells.Item(500, 10).Value = 195
$ws.Cells.Item(500, 11).Value = 11000
$ws.Cells.Item(500, 12).Value = 12000
$ws.Cells.Item(500, 13).Value = 11462
$ws.Cells.Item(500, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(500, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(500, 16).Value = 229
$ws.Cells.Item(500, 17).Value = 50
$ws.Cells.Item(500, 18).Value = "Hortaliza"
